$d = $word.ActiveDocument

# 1) Remove the existing "_GoBack" bookmark at the end of the document
#    (it currently sits after the "...symmetrical distribution." run).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) Insert a new "PTO" paragraph (with a relocated "_GoBack" bookmark)
#    plus a trailing blank paragraph, right after the empty paragraph
#    that follows the "...as the variance increases." / chart picture
#    block (paragraph 16), and before the bold heading block that starts
#    the "68-95-99.7 Rule" section (paragraph 17).
$anchor = $d.Paragraphs(16)
$anchor.Range.InsertParagraphAfter()
$anchor = $d.Paragraphs(16)
$anchor.Range.InsertParagraphAfter()

$ptoPara = $d.Paragraphs(17)
$ptoPara.Range.Text = "PTO"

$ptoPara = $d.Paragraphs(17)
$ptoRange = $ptoPara.Range
$bookmarkRange = $d.Range($ptoRange.Start, $ptoRange.End - 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
